$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "data" table (Table3) currently spans A1:E72; add the 23 May 2020
# record as a new row, growing the table to A1:E73.
$table = $ws.ListObjects.Item("Table3")
$newListRow = $table.ListRows.Add()
$newRow = $newListRow.Range.Row

$ws.Cells.Item($newRow, 1).Value = 43974   # date  -> 2020-05-23
$ws.Cells.Item($newRow, 2).Value = 40178   # test
$ws.Cells.Item($newRow, 3).Value = 1186    # case
$ws.Cells.Item($newRow, 4).Value = 32      # death
$ws.Cells.Item($newRow, 5).Value = 1491    # recovered

# Match the resulting view/selection state from the diff
$ws.Range("C73").Select()
